# "Generate Report for Handback"
# Updates the localization-status workbook with handback info for the
# e9bd2ce1-68af-4a8e-a04a-1fea47d60e61 file on both the zh-cn and de-de
# report sheets, and widens the "Error Detail" column on those sheets.

$wb = $excel.ActiveWorkbook

# Shared error message referenced from both language sheets (row 8 / P column)
$errorDetail = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/876b4ddb415b436f1cba55ad8c4780884e1fb077/e2e/e9bd2ce1-68af-4a8e-a04a-1fea47d60e61.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/816d7cba27664ec36330fb540aef01771b447898/e2e/e9bd2ce1-68af-4a8e-a04a-1fea47d60e61.md.'

# column width (stored internal width ends up as 40 when ColumnWidth is set
# to this value in this runtime's unit conversion)
$targetColumnWidth = 39.166666666666664

function Update-HandbackSheet {
    param(
        [string]$SheetName,
        [string]$TargetFileDisplay,
        [string]$TargetFileUrl,
        [string]$HandbackFileName,
        [string]$HandbackDateTime,
        [string]$ErrorDetail
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Widen the "Error Detail" column (column 16 / P) so the long message is readable
    $ws.Columns.Item(16).ColumnWidth = $targetColumnWidth

    # I8 - Latest Target File: now resolved, becomes a hyperlink to the target repo file
    $ws.Range("I8").Value = $TargetFileDisplay
    $ws.Hyperlinks.Add($ws.Range("I8"), $TargetFileUrl, "", "", $TargetFileDisplay)

    # J8 - Latest Handback File: the generated xlf handback file name
    $ws.Range("J8").Value = $HandbackFileName

    # K8 - Latest Handback DateTime
    $ws.Range("K8").Value = $HandbackDateTime

    # P8 - Error Detail
    $ws.Range("P8").Value = $ErrorDetail
}

Update-HandbackSheet "zh-cn" `
    "e9bd2ce1-68af-4a8e-a04a-1fea47d60e61.md" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/554c466824db739dca9891ca660b2e3aae029b40/e2e/e9bd2ce1-68af-4a8e-a04a-1fea47d60e61.md" `
    "e9bd2ce1-68af-4a8e-a04a-1fea47d60e61.c896026ece00bfd22729ce087b576c68e1f4d294.zh-cn.xlf" `
    "2016-09-05 06:48:40" `
    $errorDetail

Update-HandbackSheet "de-de" `
    "e9bd2ce1-68af-4a8e-a04a-1fea47d60e61.md" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4608341882d4665ebaa643b24627989d4e61ff0b/e2e/e9bd2ce1-68af-4a8e-a04a-1fea47d60e61.md" `
    "e9bd2ce1-68af-4a8e-a04a-1fea47d60e61.c896026ece00bfd22729ce087b576c68e1f4d294.de-de.xlf" `
    "2016-09-05 06:48:47" `
    $errorDetail

Write-Output "Report generated for handback."
